$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.79520333333333
$ws.Range("H2").Value = 38.38560999999999
$ws.Range("I2").Value = 0.4115337443623667
$ws.Range("J2").Value = 0.4115337443623667
$ws.Range("M2").Value = 0.66094
$ws.Range("N2").Value = 1.98282
$ws.Range("O2").Value = 0.1469991764185096
$ws.Range("P2").Value = 0.1469991764185096
$ws.Range("Q2").Value = 8.456861691133332
$ws.Range("R2").Value = 76.11175522019998
$ws.Range("S2").Value = 0.06049512148969336
$ws.Range("T2").Value = 0.06049512148969336
$ws.Range("G3").Value = 12.79520333333333
$ws.Range("H3").Value = 38.38560999999999
$ws.Range("I3").Value = 0.4115337443623667
$ws.Range("J3").Value = 0.4115337443623667
$ws.Range("M3").Value = 2.250698333333333
$ws.Range("N3").Value = 6.752095
$ws.Range("O3").Value = 0.5005761511884772
$ws.Range("P3").Value = 0.5005761511884772
$ws.Range("Q3").Value = 28.79814281699444
$ws.Range("R3").Value = 259.18328535295
$ws.Range("S3").Value = 0.2060039778370962
$ws.Range("T3").Value = 0.2060039778370962
$ws.Range("G4").Value = 12.79520333333333
$ws.Range("H4").Value = 38.38560999999999
$ws.Range("I4").Value = 0.4115337443623667
$ws.Range("J4").Value = 0.4115337443623667
$ws.Range("M4").Value = 1.584577333333334
$ws.Range("N4").Value = 4.753732
$ws.Range("O4").Value = 0.3524246723930132
$ws.Range("P4").Value = 0.3524246723930132
$ws.Range("Q4").Value = 20.27498917739111
$ws.Range("R4").Value = 182.47490259652
$ws.Range("S4").Value = 0.1450346450355771
$ws.Range("T4").Value = 0.1450346450355771
$ws.Range("I5").Value = 0.3715860904941301
$ws.Range("J5").Value = 0.3715860904941301
$ws.Range("M5").Value = 0.66094
$ws.Range("N5").Value = 1.98282
$ws.Range("O5").Value = 0.1469991764185096
$ws.Range("P5").Value = 0.1469991764185096
$ws.Range("Q5").Value = 7.635952620426666
$ws.Range("R5").Value = 68.72357358383999
$ws.Range("S5").Value = 0.05462284927121089
$ws.Range("T5").Value = 0.05462284927121089
$ws.Range("I6").Value = 0.3715860904941301
$ws.Range("J6").Value = 0.3715860904941301
$ws.Range("M6").Value = 2.250698333333333
$ws.Range("N6").Value = 6.752095
$ws.Range("O6").Value = 0.5005761511884772
$ws.Range("P6").Value = 0.5005761511884772
$ws.Range("Q6").Value = 26.00270196418222
$ws.Range("S6").Value = 0.1860071350147248
$ws.Range("T6").Value = 0.1860071350147248
$ws.Range("I7").Value = 0.3715860904941301
$ws.Range("J7").Value = 0.3715860904941301
$ws.Range("M7").Value = 1.584577333333334
$ws.Range("N7").Value = 4.753732
$ws.Range("O7").Value = 0.3524246723930132
$ws.Range("P7").Value = 0.3524246723930132
$ws.Range("Q7").Value = 18.30689236653156
$ws.Range("R7").Value = 164.762031298784
$ws.Range("S7").Value = 0.1309561062081943
$ws.Range("T7").Value = 0.1309561062081943
$ws.Range("G8").Value = 6.743130666666667
$ws.Range("H8").Value = 20.229392
$ws.Range("I8").Value = 0.2168801651435032
$ws.Range("J8").Value = 0.2168801651435032
$ws.Range("M8").Value = 0.66094
$ws.Range("N8").Value = 1.98282
$ws.Range("O8").Value = 0.1469991764185096
$ws.Range("P8").Value = 0.1469991764185096
$ws.Range("Q8").Value = 4.456804782826667
$ws.Range("R8").Value = 40.11124304544
$ws.Range("S8").Value = 0.03188120565760531
$ws.Range("T8").Value = 0.03188120565760531
$ws.Range("G9").Value = 6.743130666666667
$ws.Range("H9").Value = 20.229392
$ws.Range("I9").Value = 0.2168801651435032
$ws.Range("J9").Value = 0.2168801651435032
$ws.Range("M9").Value = 2.250698333333333
$ws.Range("N9").Value = 6.752095
$ws.Range("O9").Value = 0.5005761511884772
$ws.Range("P9").Value = 0.5005761511884772
$ws.Range("Q9").Value = 15.17675295291555
$ws.Range("R9").Value = 136.59077657624
$ws.Range("S9").Value = 0.1085650383366562
$ws.Range("T9").Value = 0.1085650383366562
$ws.Range("G10").Value = 6.743130666666667
$ws.Range("H10").Value = 20.229392
$ws.Range("I10").Value = 0.2168801651435032
$ws.Range("J10").Value = 0.2168801651435032
$ws.Range("M10").Value = 1.584577333333334
$ws.Range("N10").Value = 4.753732
$ws.Range("O10").Value = 0.3524246723930132
$ws.Range("P10").Value = 0.3524246723930132
$ws.Range("Q10").Value = 10.68501201010489
$ws.Range("R10").Value = 96.165108090944
$ws.Range("S10").Value = 0.0764339211492417
$ws.Range("T10").Value = 0.0764339211492417
